$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 from the numeric date value to the new shared string
$ws.Range("A2").Value = "BC.ICFKHI0000001"

# Update the active selection to F14
$ws.Range("F14").Select()
